# Update the cached "datetimeFigureOut" date field text from 1/21/2019 to
# 3/10/2020 wherever it appears: the Slide Master's Date Placeholder and
# every Slide Layout's Date Placeholder.

$p = $ppt.ActivePresentation

$oldDate = "1/21/2019"
$newDate = "3/10/2020"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
